$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B2").Value = 17.812
$ws.Range("C2").Value = 12.66
$ws.Range("D2").Value = 29.056
$ws.Range("E2").Value = 0.8247863247863249
$ws.Range("F2").Value = 0.3381555153707053

$ws.Range("G4").Value = 0.3968835930339138
$ws.Range("G5").Value = 0.7687901008249315
$ws.Range("G6").Value = 0.5680568285976172
$ws.Range("G7").Value = 0.4718148487626035
